$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rename the city "Caldas" -> "Manizales" (row 8 / column H of the matrix).
# Updating the two label cells updates the shared-string table automatically.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Manizales"
$ws.Range("H1").Value = "Manizales"

# ---------------------------------------------------------------------------
# New travel times between Manizales and the other cities.
# Column H (rows 2-16) and row 8 (columns B-P) both represent the
# "Manizales" column/row of the symmetric travel-time matrix.
# ---------------------------------------------------------------------------
$hValues = @(1.97, 16.72, 5.83, 8.77, 15.27, 13.45, 0, 3.72, 11.56, 3.95, 1.17, 7.34, 16.72, 3.95, 13.32)

for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $hValues[$i]   # column H
}

$rowCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")
for ($i = 0; $i -lt $rowCols.Length; $i++) {
    $ws.Range($rowCols[$i] + "8").Value = $hValues[$i]   # row 8
}

# ---------------------------------------------------------------------------
# Column I (rows 2-16) and row 9 (columns B-P) both represent the
# "Medellin" column/row; these previously held formulas that referenced the
# Manizales (Caldas) column/row. They are now plain (computed) values.
# ---------------------------------------------------------------------------
$iValues = @(5.25, 12.03, 8.17, 8.1999999999999993, 11.8, 11.38, 3.72, 0, 7.42, 7.75, 4.08, 14.48, 12.03, 7.75, 11.55)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]   # column I
}

for ($i = 0; $i -lt $rowCols.Length; $i++) {
    $ws.Range($rowCols[$i] + "9").Value = $iValues[$i]   # row 9
}

# ---------------------------------------------------------------------------
# Restore the selection/view state recorded after the edit.
# ---------------------------------------------------------------------------
$ws.Range("P9").Select()
$excel.ActiveWindow.ScrollRow = 1
